$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 300, pushing existing rows 300-313 down to 302-315.
$ws.Range("A300:R301").EntireRow.Insert()

# New row 300 data
$ws.Cells.Item(300, 1).Value = 5
$ws.Cells.Item(300, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(300, 3).Value = "Maule"
$ws.Cells.Item(300, 4).Value = 44509
$ws.Cells.Item(300, 5).Value = 7
$ws.Cells.Item(300, 6).Value = 100114001
$ws.Cells.Item(300, 7).Value = "Papa"
$ws.Cells.Item(300, 8).Value = "Rodeo"
$ws.Cells.Item(300, 9).Value = "1a nueva(o)"
$ws.Cells.Item(300, 10).Value = 1200
$ws.Cells.Item(300, 11).Value = 9000
$ws.Cells.Item(300, 12).Value = 9000
$ws.Cells.Item(300, 13).Value = 9000
$ws.Cells.Item(300, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(300, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(300, 16).Value = 360
$ws.Cells.Item(300, 17).Value = 25
$ws.Cells.Item(300, 18).Value = "Hortaliza"

# New row 301 data
$ws.Cells.Item(301, 1).Value = 5
$ws.Cells.Item(301, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(301, 3).Value = "Maule"
$ws.Cells.Item(301, 4).Value = 44509
$ws.Cells.Item(301, 5).Value = 7
$ws.Cells.Item(301, 6).Value = 100114001
$ws.Cells.Item(301, 7).Value = "Papa"
$ws.Cells.Item(301, 8).Value = "Rosara"
$ws.Cells.Item(301, 9).Value = "1a nueva(o)"
$ws.Cells.Item(301, 10).Value = 1600
$ws.Cells.Item(301, 11).Value = 9000
$ws.Cells.Item(301, 12).Value = 9000
$ws.Cells.Item(301, 13).Value = 9000
$ws.Cells.Item(301, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(301, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(301, 16).Value = 360
$ws.Cells.Item(301, 17).Value = 25
$ws.Cells.Item(301, 18).Value = "Hortaliza"

# Apply the date format style (same as other D column cells) to the new D cells
$ws.Range("D300:D301").NumberFormat = $ws.Range("D302").NumberFormat
